$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "SUMMARY OF DISCUSSION" heading: the text was split across two
#    runs ("summary of Discus" + "sion") with a _GoBack bookmark sitting
#    between them. Re-writing the whole phrase in one shot merges the
#    runs back into a single run and drops the now-crossed bookmark.
# ------------------------------------------------------------------
$d.Content.Find.Execute("summary of Discussion", $false, $false, $false, `
    $false, $false, $true, 1, $false, "summary of Discussion", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Rewrite of the agenda paragraph. Apply the wording changes in
#    left-to-right order so each Find starts from a still-unique
#    anchor of surrounding, unmodified text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("to brainstorm some ideas and to select one for", $false, $false, $false, `
    $false, $false, $true, 1, $false, "to brainstorm for ideas and inspiration, and select a realistic one for", 2) | Out-Null

$d.Content.Find.Execute("We decided to go with android studio developing an android chat application which will work", $false, $false, $false, `
    $false, $false, $true, 1, $false, "We decided to use android studio as our IDE to develop an android chat application that will work", 2) | Out-Null

$d.Content.Find.Execute("Every group member was asked to learn socket programming.", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Every group member was asked to learn about socket programming.", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark at the new last-edit location:
#    right before "socket programming." (after "... learn about ").
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Text = "socket programming"
$r.Find.Execute() | Out-Null
$goBack = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null

# ------------------------------------------------------------------
# 4) "Akhtar" / "Zaman" were two separate (spell-check-flagged) runs
#    with a space run between them; merge into a single "Akhtar Zaman"
#    run just like the heading fix above.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Akhtar Zaman", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Akhtar Zaman", 2) | Out-Null
